# Calibration Curves Generated for Pressure
#
# Adds per-column average-calibration labels/values (avg1..avg5) next to
# the raw calibration data on the "Simple Data" sheet, mirroring the
# author's manual entry of J2:K6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Data")

# Label + average formula pairs, one per calibration channel (D..H).
# Row 2's "average" was entered as a bare range reference (D2:D66),
# which Excel reduces to the first cell of the range (1462) rather than
# a true AVERAGE() — preserved here verbatim to match the source data.
$ws.Range("J2").Value = "avg1"
$ws.Range("K2").Formula = "=D2:D66"

$ws.Range("J3").Value = "avg2"
$ws.Range("K3").Formula = "=AVERAGE(E2:E66)"

$ws.Range("J4").Value = "avg3"
$ws.Range("K4").Formula = "=AVERAGE(F2:F66)"

$ws.Range("J5").Value = "avg4"
$ws.Range("K5").Formula = "=AVERAGE(G2:G66)"

$ws.Range("J6").Value = "avg5"
$ws.Range("K6").Formula = "=AVERAGE(H2:H66)"

# Leave the selection where the author ended up after typing the formulas.
[void]$ws.Range("K7").Select()
